# This edit re-sorts the 4 "Hydnellum/Gomphus" observation rows (rows 6-9)
# on the "Artfynd" sheet into updated Taxonsorteringsordning (column B) order.
# The row payloads (columns A, C:AY) simply rotate between rows 6/7/8/9; only
# column B (Taxonsorteringsordning) receives brand-new values.
#
# Net effect (old row -> new row):
#   old row 8 -> new row 6   (B becomes 90806)
#   old row 9 -> new row 7   (B becomes 90814)
#   old row 7 -> new row 8   (B becomes 90832)
#   old row 6 -> new row 9   (B becomes 89057)
#
# We move the data with Range.Copy (not Range.Value=) so that every cell's
# original type/text is preserved exactly (in particular the plain-text
# "yyyy-mm-dd" strings in columns Y/AA must stay text, not get reinterpreted
# as real dates the way a fresh Value assignment would do it), and so that
# cells which were present-but-empty stay present-but-empty while cells that
# never existed in the source row stay absent in the destination row.
#
# Because this is a 4-way cyclic rotation we stage through scratch rows
# (200-203) first so no source row is overwritten before it has been copied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRows     = @(6, 7, 8, 9)
$scratchRows = @(200, 201, 202, 203)

# --- Step 1: stage current rows 6-9 into scratch rows 200-203 -------------
for ($i = 0; $i -lt $srcRows.Length; $i++) {
    $s = $srcRows[$i]
    $t = $scratchRows[$i]
    $ws.Range("A" + $s + ":AY" + $s).Copy($ws.Range("A" + $t + ":AY" + $t))
}

# --- Step 2: copy the scratch rows back into their new homes --------------
# new row 6 <- old row 8 (scratch 202)
$ws.Range("A202:AY202").Copy($ws.Range("A6:AY6"))
# new row 7 <- old row 9 (scratch 203)
$ws.Range("A203:AY203").Copy($ws.Range("A7:AY7"))
# new row 8 <- old row 7 (scratch 201)
$ws.Range("A201:AY201").Copy($ws.Range("A8:AY8"))
# new row 9 <- old row 6 (scratch 200)
$ws.Range("A200:AY200").Copy($ws.Range("A9:AY9"))

# --- Step 3: whole-row copy also drags along empty cells for columns that -
# never actually appear in rows 6-9 (e.g. L, M, O, ...); strip those back
# out so the destination rows only contain the columns that really belong.
$neverUsedRanges = @("L6:M9", "O6:O9", "X6:X9", "AC6:AC9", "AJ6:AS9", "AU6:AV9")
foreach ($rng in $neverUsedRanges) {
    $ws.Range($rng).ClearContents()
}

# --- Step 4: write the updated Taxonsorteringsordning (column B) values ---
$ws.Range("B6").Value = 90806
$ws.Range("B7").Value = 90814
$ws.Range("B8").Value = 90832
$ws.Range("B9").Value = 89057

# --- Step 5: clean up the scratch rows -------------------------------------
$ws.Range("A200:AY203").ClearContents()

Write-Output "Row re-sort applied."
